# Auto-generated edit script applying scheduled market-data refresh to Sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 298.4
$ws.Range("I2").Value = 190.66667
$ws.Range("J2").Value = 460
$ws.Range("K2").Value = 190.66667
$ws.Range("L2").Value = 460
$ws.Range("M2").Value = -77.66667000000001
$ws.Range("N2").Value = -686

$ws.Range("H19").Value = 758.0769
$ws.Range("J19").Value = 955.4
$ws.Range("L19").Value = 955.4
$ws.Range("N19").Value = -1305.4

$ws.Range("H28").Value = 3347
$ws.Range("I28").Value = 1218.7
$ws.Range("J28").Value = 5120.5835
$ws.Range("K28").Value = 1218.7
$ws.Range("L28").Value = 5120.5835
$ws.Range("M28").Value = -733.7
$ws.Range("N28").Value = -6090.5835

$ws.Range("H88").Value = 1144.8
$ws.Range("I88").Value = 900
$ws.Range("J88").Value = 1206
$ws.Range("K88").Value = 900
$ws.Range("L88").Value = 1206
$ws.Range("M88").Value = -494
$ws.Range("N88").Value = -2018

$ws.Range("H91").Value = 1144.8
$ws.Range("I91").Value = 900
$ws.Range("J91").Value = 1206
$ws.Range("K91").Value = 900
$ws.Range("L91").Value = 1206
$ws.Range("M91").Value = 504
$ws.Range("N91").Value = -4014

$ws.Range("H107").Value = 1910
$ws.Range("I107").Value = 1910
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1910
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 10
$ws.Range("N107").Value = ""

$ws.Range("H112").Value = 1674.7142
$ws.Range("J112").Value = 1294.9
$ws.Range("L112").Value = 3884.7
$ws.Range("N112").Value = -6100.700000000001

$ws.Range("H113").Value = 8259.454
$ws.Range("I113").Value = 8111.5557
$ws.Range("K113").Value = 8111.5557
$ws.Range("M113").Value = -4857.5557

$ws.Range("H131").Value = 1606.75
$ws.Range("I131").Value = 1606.75
$ws.Range("K131").Value = 4820.25
$ws.Range("M131").Value = 219.75

$ws.Range("H132").Value = 1953.0834
$ws.Range("I132").Value = 1472.7826
$ws.Range("K132").Value = 4418.3478
$ws.Range("M132").Value = -1888.3478

$ws.Range("H137").Value = 3521.125
$ws.Range("I137").Value = 1424
$ws.Range("K137").Value = 4272
$ws.Range("M137").Value = -1722

$ws.Range("H138").Value = 3919.3333
$ws.Range("I138").Value = 1766.3334
$ws.Range("J138").Value = 6072.3335
$ws.Range("K138").Value = 5299.0002
$ws.Range("L138").Value = 18217.0005
$ws.Range("M138").Value = -159.0002000000004
$ws.Range("N138").Value = -28497.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2825.0557
$ws.Range("I45").Value = 2051.5833
$ws.Range("K45").Value = 2051.5833
$ws.Range("M45").Value = -1674.5833

$ws.Range("H74").Value = 4726.25
$ws.Range("I74").Value = 4309.1
$ws.Range("K74").Value = 4309.1
$ws.Range("M74").Value = -3435.1

$ws.Range("H77").Value = 4726.25
$ws.Range("I77").Value = 4309.1
$ws.Range("K77").Value = 21545.5
$ws.Range("M77").Value = -17177.5

$ws.Range("H88").Value = 1150
$ws.Range("I88").Value = 1107
$ws.Range("J88").Value = 1180.7142
$ws.Range("K88").Value = 1107
$ws.Range("L88").Value = 1180.7142
$ws.Range("M88").Value = -701
$ws.Range("N88").Value = -1992.7142

$ws.Range("H91").Value = 1150
$ws.Range("I91").Value = 1107
$ws.Range("J91").Value = 1180.7142
$ws.Range("K91").Value = 1107
$ws.Range("L91").Value = 1180.7142
$ws.Range("M91").Value = 297
$ws.Range("N91").Value = -3988.7142

$ws.Range("H122").Value = 2994
$ws.Range("I122").Value = 2994
$ws.Range("K122").Value = 8982
$ws.Range("M122").Value = -6532

$ws.Range("H126").Value = 5250
$ws.Range("I126").Value = 5250
$ws.Range("K126").Value = 15750
$ws.Range("M126").Value = -13280

$ws.Range("H132").Value = 1618.4762
$ws.Range("I132").Value = 1317
$ws.Range("K132").Value = 3951
$ws.Range("M132").Value = -1421

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 556.6
$ws.Range("I80").Value = 300.125
$ws.Range("K80").Value = 300.125
$ws.Range("M80").Value = 697.875

$ws.Range("H83").Value = 556.6
$ws.Range("I83").Value = 300.125
$ws.Range("K83").Value = 1500.625
$ws.Range("M83").Value = 3491.375

$ws.Range("H105").Value = 6994590.5
$ws.Range("I105").Value = 6994590.5
$ws.Range("K105").Value = 6994590.5
$ws.Range("M105").Value = -6992843.5

$ws.Range("H107").Value = 45459492
$ws.Range("I107").Value = 125000584
$ws.Range("K107").Value = 125000584
$ws.Range("M107").Value = -124998664

$ws.Range("H134").Value = 1763.875
$ws.Range("I134").Value = 1273.85
$ws.Range("J134").Value = 4214
$ws.Range("K134").Value = 3821.55
$ws.Range("L134").Value = 12642
$ws.Range("M134").Value = -1286.55
$ws.Range("N134").Value = -17712

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4524.769
$ws.Range("I31").Value = 1476.2667
$ws.Range("K31").Value = 1476.2667
$ws.Range("M31").Value = -1181.2667

$ws.Range("H34").Value = 4524.769
$ws.Range("I34").Value = 1476.2667
$ws.Range("K34").Value = 1476.2667
$ws.Range("M34").Value = -1274.2667

$ws.Range("H107").Value = 1403.7333
$ws.Range("I107").Value = 512.25
$ws.Range("J107").Value = 2422.5715
$ws.Range("K107").Value = 512.25
$ws.Range("L107").Value = 2422.5715
$ws.Range("M107").Value = 1407.75
$ws.Range("N107").Value = -6262.5715

$ws.Range("H134").Value = 2379.6667
$ws.Range("I134").Value = 2379.6667
$ws.Range("K134").Value = 7139.000100000001
$ws.Range("M134").Value = -4604.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 554.4737
$ws.Range("I107").Value = 265
$ws.Range("J107").Value = 765
$ws.Range("K107").Value = 795
$ws.Range("L107").Value = 2295
$ws.Range("M107").Value = 1125
$ws.Range("N107").Value = -6135

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 533.8
$ws.Range("I97").Value = 292.25
$ws.Range("K97").Value = 292.25
$ws.Range("M97").Value = 203.75

$ws.Range("H102").Value = 3157.5625
$ws.Range("I102").Value = 2402.2
$ws.Range("J102").Value = 4416.5
$ws.Range("K102").Value = 2402.2
$ws.Range("L102").Value = 4416.5
$ws.Range("M102").Value = -780.1999999999998
$ws.Range("N102").Value = -7660.5

$ws.Range("H107").Value = 718.5
$ws.Range("I107").Value = 167
$ws.Range("K107").Value = 167
$ws.Range("M107").Value = 1753

$ws.Range("H122").Value = 1965.7142
$ws.Range("I122").Value = 1131.8334
$ws.Range("J122").Value = 6969
$ws.Range("K122").Value = 3395.5002
$ws.Range("L122").Value = 20907
$ws.Range("M122").Value = -945.5001999999999
$ws.Range("N122").Value = -25807

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 937.35297
$ws.Range("J16").Value = 239
$ws.Range("L16").Value = 239
$ws.Range("N16").Value = -579

$ws.Range("H32").Value = 1008.6667
$ws.Range("I32").Value = 1008.6667
$ws.Range("K32").Value = 1008.6667
$ws.Range("M32").Value = -691.6667

$ws.Range("H61").Value = 125005820
$ws.Range("I61").Value = 200004510
$ws.Range("K61").Value = 200004510
$ws.Range("M61").Value = -200004308

$ws.Range("H113").Value = 125005820
$ws.Range("I113").Value = 200004510
$ws.Range("K113").Value = 200004510
$ws.Range("M113").Value = -200002340

$ws.Range("H122").Value = 3495.7144
$ws.Range("I122").Value = 3495.7144
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10487.1432
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8037.143199999999
$ws.Range("N122").Value = ""

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

$ws.Range("H134").Value = 25000
$ws.Range("J134").Value = 25000
$ws.Range("L134").Value = 25000
$ws.Range("N134").Value = -35140

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 27500
$ws.Range("J16").Value = 35000
$ws.Range("L16").Value = 35000
$ws.Range("N16").Value = -35584

$ws.Range("H133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
